# Updates cryptocurrency Price (column D) and Volume(1h) % change (column E)
# values in the cryptos list, refreshed by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '37.042.48'
$ws.Cells.Item(2, 5).Value = '  -0.34%  '
$ws.Cells.Item(3, 4).Value = '2.044.05'
$ws.Cells.Item(3, 5).Value = '  -0.66%  '
$ws.Cells.Item(4, 5).Value = '  +0.08%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '246.05'
$ws.Cells.Item(5, 5).Value = '  -1.71%  '
$ws.Cells.Item(6, 5).Value = '  -2.02%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '58.94'
$ws.Cells.Item(7, 5).Value = '  -3.22%  '
$ws.Cells.Item(8, 5).Value = '  +0.00%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.378'
$ws.Cells.Item(9, 5).Value = '  -2.16%  '
$ws.Cells.Item(10, 5).Value = '  -2.37%  '
$ws.Cells.Item(11, 5).Value = '  +2.20%  '
$ws.Cells.Item(12, 5).Value = '  -5.39%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.893'
$ws.Cells.Item(13, 5).Value = '  +8.05%  '
$ws.Cells.Item(14, 4).Value = '2.341.67'
$ws.Cells.Item(14, 5).Value = '  -0.63%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '5.71'
$ws.Cells.Item(15, 5).Value = '  -0.14%  '
$ws.Cells.Item(16, 4).Value = '2.039.29'
$ws.Cells.Item(16, 5).Value = '  -0.86%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '18.44'
$ws.Cells.Item(17, 5).Value = '  +2.25%  '
$ws.Cells.Item(18, 4).Value = '37.023.63'
$ws.Cells.Item(18, 5).Value = '  -0.42%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '73.66'
$ws.Cells.Item(19, 5).Value = '  -2.25%  '
$ws.Cells.Item(20, 5).Value = '  -2.18%  '
$ws.Cells.Item(21, 5).Value = '  -0.80%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '240.37'
$ws.Cells.Item(22, 5).Value = '  +0.28%  '
$ws.Cells.Item(23, 5).Value = '  -0.05%  '
$ws.Cells.Item(24, 5).Value = '  +1.64%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '9.66'
$ws.Cells.Item(25, 5).Value = '  +2.34%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '168.37'
$ws.Cells.Item(26, 5).Value = '  -0.63%  '
$ws.Cells.Item(27, 5).Value = '  -3.73%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '19.98'
$ws.Cells.Item(28, 5).Value = '  -0.34%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '5.53'
$ws.Cells.Item(29, 5).Value = '  +14.67%  '
$ws.Cells.Item(30, 5).Value = '  -1.08%  '
$ws.Cells.Item(31, 5).Value = '  -2.71%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.75'
$ws.Cells.Item(32, 5).Value = '  +4.28%  '
$ws.Cells.Item(33, 5).Value = '  -1.31%  '
$ws.Cells.Item(34, 5).Value = '  +0.24%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.84'
$ws.Cells.Item(35, 5).Value = '  +6.06%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.0850'
$ws.Cells.Item(36, 5).Value = '  -5.30%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.25'
$ws.Cells.Item(37, 5).Value = '  -0.76%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.31'
$ws.Cells.Item(38, 5).Value = '  -4.01%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.25'
$ws.Cells.Item(39, 5).Value = '  -1.63%  '
$ws.Cells.Item(40, 5).Value = '  -2.18%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0978'
$ws.Cells.Item(41, 5).Value = '  -10.26%  '
$ws.Cells.Item(42, 5).Value = '  -0.65%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.14'
$ws.Cells.Item(43, 5).Value = '  +0.33%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '97.36'
$ws.Cells.Item(44, 5).Value = '  -0.71%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '17.00'
$ws.Cells.Item(45, 5).Value = '  -6.70%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.38'
$ws.Cells.Item(46, 5).Value = '  -4.49%  '
$ws.Cells.Item(47, 4).Value = '1.298.48'
$ws.Cells.Item(47, 5).Value = '  +0.13%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.86'
$ws.Cells.Item(48, 5).Value = '  -0.28%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '6.74'
$ws.Cells.Item(49, 5).Value = '  -1.77%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '3.66'
$ws.Cells.Item(50, 5).Value = '  +2.53%  '
$ws.Cells.Item(51, 4).Value = '2.227.07'
$ws.Cells.Item(51, 5).Value = '  -0.56%  '
